$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Auto-update data + news: refresh Job Openings (JTSJOL_mil) row with latest release
$ws.Range("E7").Value = 6.542

# Force F7 to remain plain text (avoid Excel auto-converting "Dec 2025" into a date serial)
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "Dec 2025"
$ws.Range("F7").ClearFormats()

$ws.Range("G7").Value = 7.796471074380164
$ws.Range("H7").Value = -0.9660000000000002
$ws.Range("I7").Value = -0.1286627597229622
